$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'98.219.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.30%  "

# Row 3
$ws.Range("D3").Value = "'3.405.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.15%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'255.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.06%  "

# Row 6
$ws.Range("D6").Value = "'680.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.66%  "

# Row 7
$ws.Range("D7").Value = "'1.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.80%  "

# Row 8
$ws.Range("D8").Value = "'0.431"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.91%  "

# Row 9
$ws.Range("E9").Value = "  -4.44%  "

# Row 10
$ws.Range("E10").Value = "  -0.01%  "

# Row 11
$ws.Range("D11").Value = "'3.402.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.20%  "

# Row 12
$ws.Range("E12").Value = "  +2.04%  "

# Row 13
$ws.Range("D13").Value = "'41.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.70%  "

# Row 14
$ws.Range("D14").Value = "'6.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +11.79%  "

# Row 15
$ws.Range("D15").Value = "'97.740.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.12%  "

# Row 16
$ws.Range("D16").Value = "'0.0000266"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.20%  "

# Row 17
$ws.Range("D17").Value = "'4.038.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.13%  "

# Row 18
$ws.Range("D18").Value = "'8.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +16.11%  "

# Row 19
$ws.Range("D19").Value = "'3.404.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.08%  "

# Row 20
$ws.Range("D20").Value = "'0.576"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +28.39%  "

# Row 21
$ws.Range("D21").Value = "'17.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.46%  "

# Row 22
$ws.Range("D22").Value = "'10.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.24%  "

# Row 23
$ws.Range("D23").Value = "'3.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.59%  "

# Row 24
$ws.Range("D24").Value = "'508.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.18%  "

# Row 25
$ws.Range("D25").Value = "'0.0000204"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.57%  "

# Row 26
$ws.Range("D26").Value = "'6.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.04%  "

# Row 27
$ws.Range("D27").Value = "'100.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.27%  "

# Row 28
$ws.Range("D28").Value = "'12.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "

# Row 29
$ws.Range("D29").Value = "'3.586.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.03%  "

# Row 30
$ws.Range("D30").Value = "'0.151"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.41%  "

# Row 31
$ws.Range("D31").Value = "'11.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.34%  "

# Row 32
$ws.Range("E32").Value = "  +0.02%  "

# Row 33
$ws.Range("D33").Value = "'0.195"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.82%  "

# Row 34
$ws.Range("E34").Value = "  +22.82%  "

# Row 35
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
$ws.Range("D36").Value = "'0.572"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.61%  "

# Row 37
$ws.Range("D37").Value = "'29.65"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.33%  "

# Row 38
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'1.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.24%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'7.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.68%  "

# Row 40
$ws.Range("D40").Value = "'531.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.17%  "

# Row 41
$ws.Range("D41").Value = "'0.153"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.99%  "

# Row 42
$ws.Range("E42").Value = "  -0.03%  "

# Row 43
$ws.Range("D43").Value = "'24.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.12%  "

# Row 44
$ws.Range("D44").Value = "'0.869"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.44%  "

# Row 45
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0439"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.14%  "

# Row 46
$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").Value = "'9.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.06%  "

# Row 47
$ws.Range("B47").Value = "MantraDAO"
$ws.Range("C47").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D47").Value = "'3.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.95%  "

# Row 48
$ws.Range("D48").Value = "'1.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +14.03%  "

# Row 49
$ws.Range("D49").Value = "'5.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.18%  "

# Row 50
$ws.Range("D50").Value = "'56.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.14%  "

# Row 51
$ws.Range("D51").Value = "'3.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.93%  "

